# Apply updated crypto price/volume data as per Sat Oct  5 10:39:06 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.254.73"
$ws.Range("E2").Value = "  +1.54%  "

$ws.Range("D3").Value = "2.423.02"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("D5").Value = "'563.31"
$ws.Range("E5").Value = "  +2.11%  "

$ws.Range("D6").Value = "'144.26"
$ws.Range("E6").Value = "  +3.26%  "

$ws.Range("E8").Value = "  +1.67%  "

$ws.Range("D9").Value = "2.421.70"

$ws.Range("D10").Value = "'0.109"
$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("E11").Value = "  -2.11%  "

$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("E13").Value = "  +0.61%  "

$ws.Range("D14").Value = "'25.93"
$ws.Range("E14").Value = "  +1.77%  "

$ws.Range("E15").Value = "  +5.15%  "

$ws.Range("D16").Value = "2.861.45"
$ws.Range("E16").Value = "  +2.01%  "

$ws.Range("D17").Value = "62.091.72"
$ws.Range("E17").Value = "  +1.22%  "

$ws.Range("D18").Value = "2.420.34"
$ws.Range("E18").Value = "  +1.71%  "

$ws.Range("D19").Value = "'11.33"
$ws.Range("E19").Value = "  +3.19%  "

$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").Value = "'323.99"
$ws.Range("E21").Value = "  +1.06%  "

$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "'65.60"
$ws.Range("E24").Value = "  +1.97%  "

$ws.Range("D25").Value = "'1.70"
$ws.Range("E25").Value = "  -2.95%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.530.77"
$ws.Range("E28").Value = "  +1.53%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").Value = "0.0₃0943"
$ws.Range("E30").Value = "  +4.88%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.27"
$ws.Range("E31").Value = "  +0.89%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.46"
$ws.Range("E32").Value = "  +4.90%  "

$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("E34").Value = "  +2.04%  "

$ws.Range("E35").Value = "  +1.57%  "

$ws.Range("D36").Value = "'5.71"
$ws.Range("E36").Value = "  +3.69%  "

$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("E38").Value = "  +2.37%  "

$ws.Range("E39").Value = "  +1.60%  "

$ws.Range("D40").Value = "'153.13"
$ws.Range("E40").Value = "  +4.25%  "

$ws.Range("D41").Value = "'18.66"
$ws.Range("E41").Value = "  +0.93%  "

$ws.Range("D42").Value = "'1.82"
$ws.Range("E42").Value = "  -3.28%  "

$ws.Range("D44").Value = "'2.32"
$ws.Range("E44").Value = "  +7.57%  "

$ws.Range("D45").Value = "'150.19"
$ws.Range("E45").Value = "  +1.61%  "

$ws.Range("D46").Value = "'3.66"
$ws.Range("E46").Value = "  +1.54%  "

$ws.Range("E47").Value = "  +2.64%  "

$ws.Range("D48").Value = "'20.34"
$ws.Range("E48").Value = "  +3.64%  "

$ws.Range("E49").Value = "  +2.31%  "

$ws.Range("D50").Value = "'0.0925"
$ws.Range("E50").Value = "  +2.14%  "

$ws.Range("E51").Value = "  +1.95%  "
